$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 359

# Insert new row 3 (0 / 344), shifting old row3 (1/44) down - but since diff shows
# final values explicitly, just set all target cells directly.
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 344

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 109

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 97

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 92

# Copy the formatting of A2 (existing styled cell) down to the newly added A4:A6 cells
# so they match the bordered/bold/centered style used by A2:A3.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4:A6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0
